# Notes for "The precedence of grouping my code with brace ()"
#
# The diff:
#  1. Removes the <w:bookmarkStart .../> for "_GoBack" from the first
#     paragraph (it was sitting right before the first run).
#  2. Splits the document so that two new (empty, centered, 36pt)
#     paragraphs are appended after the "Brace has the HIGHEST..."
#     paragraph, and the "_GoBack" bookmark (both its start and end)
#     is now wrapped around the very last (new, empty) paragraph.
#
# In other words: the cursor was left at the very end of the document
# (hence "_GoBack" moved there) after the author pressed Enter twice at
# the end of the last paragraph.

$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark completely (it will be
#    re-created at the end of the document below). "_GoBack" is a
#    hidden bookmark so it must be looked up by name rather than via
#    Bookmarks.Count/iteration.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# 2. Append two new paragraphs at the very end of the document, each
#    sharing the centered / 36pt formatting of the preceding
#    paragraph, and put a fresh "_GoBack" bookmark (empty, collapsed)
#    in the very last paragraph - exactly matching what Word leaves
#    behind after the final edit position.
$endOfDoc = $d.Content.End
$insertionPoint = $d.Range($endOfDoc, $endOfDoc)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$newParagraphsXml = "<w:p $wNs><w:pPr><w:jc w:val='center'/><w:rPr><w:sz w:val='36'/><w:szCs w:val='36'/></w:rPr></w:pPr></w:p>" + `
    "<w:p $wNs><w:pPr><w:jc w:val='center'/><w:rPr><w:sz w:val='36'/><w:szCs w:val='36'/></w:rPr></w:pPr>" + `
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/></w:p>"

[void]$insertionPoint.InsertXML($newParagraphsXml)

Write-Host "Paragraphs after edit:" $d.Paragraphs.Count
